$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The title paragraph ("BON PERMINTAAN BARANG PERSEDIAAN") gets its
#    paragraph-mark run properties (pPr/rPr) augmented with <w:b/> and
#    <w:w w:val="99"/>. The run-level formatting of the existing text is left
#    untouched. We rebuild the paragraph through InsertXML (which, for a
#    fully emptied paragraph, accepts both the pPr and the run content we
#    hand it) after temporarily clearing its text.
# ---------------------------------------------------------------------------

$titlePar = $d.Paragraphs.First
$titleRange = $titlePar.Range
$titleTextEnd = $titleRange.End - 1    # exclude the paragraph mark itself
$titleTextRange = $d.Range($titleRange.Start, $titleTextEnd)
$titleTextRange.Text = ""

$titlePar2 = $d.Paragraphs.First
$titleMarkRange = $titlePar2.Range

$titleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="14993B0D" w14:textId="3223E15D" w:rsidR="00C90568" w:rsidRPr="003164EC" w:rsidRDefault="00000000" w:rsidP="00157359"><w:pPr><w:spacing w:before="28"/><w:ind w:right="17"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:w w:val="99"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>BON</w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:spacing w:val="-4"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>PERM</w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:spacing w:val="1"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>I</w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>NTAAN</w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:spacing w:val="-13"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003164EC"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:w w:val="99"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>BARANG</w:t></w:r><w:r w:rsidR="00157359"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:w w:val="99"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> PERSEDIAAN</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$titleMarkRange.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) A new paragraph "Nomor : ${nomor}" is inserted right after the title
#    paragraph, reusing the title paragraph's pPr (same spacing/indent/
#    justification), bold + condensed (w=99) run formatting, and the
#    proofErr spell-check markers Word would normally emit around "Nomor".
# ---------------------------------------------------------------------------

$titlePar3 = $d.Paragraphs.First
$titlePar3.Range.InsertParagraphAfter()
$nomorPar = $d.Paragraphs.Item(2)
$nomorRange = $nomorPar.Range

$nomorXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:before="28"/><w:ind w:right="17"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:w w:val="99"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Nomor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsia="Arial Narrow" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:w w:val="99"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> : ${nomor}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$nomorRange.InsertXML($nomorXml)

Write-Output "edit complete"
